$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (final strike differential) column with re-pulled / recomputed values.
$updates = @{
    "F3"  = 0
    "F4"  = 5
    "F5"  = -1
    "F10" = 7
    "F13" = 6
    "F16" = -1
    "F19" = 1
    "F23" = 0
    "F24" = 1
    "F26" = 1
    "F28" = 5
    "F29" = 0
    "F30" = 1
    "F34" = 0
    "F36" = 4
    "F37" = 0
    "F47" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
